# RF010 - Gerenciar Periodos Avaliativos: 1.3 -> 1.4
# - Fixes wording of the "Data Inicial/Data Final" step + its expected result
#   (affects every test case block that reuses that shared text: TC1, TC2/4/6 steps)
# - Fixes the TC2/TC3 "confirm deletion" vs "do not confirm deletion" steps,
#   which had been swapped (and adds the missing "nao" word to the listing text)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Wording fix: 'Data Inicial e Data Final' -> 'Data Inicial' e 'Data Final' ---
$oldStep = "Lider de Pessoas preenche o campo 'Data Inicial e Data Final' informando as respectivas datas referente ao periodo"
$newStep = "Lider de Pessoas preenche o campo 'Data Inicial' e 'Data Final' informando as respectivas datas referentes ao periodo"

$oldResult = "SYSTEM apresenta o campo 'Data Inicial e Data Final' preenchido corretamente"
$newResult = "SYSTEM apresenta o campo 'Data Inicial' e 'Data Final' preenchido corretamente"

foreach ($ref in @("B14", "B46", "B65")) {
    $ws.Range($ref).Value = $newStep
}
foreach ($ref in @("D14", "D46", "D65")) {
    $ws.Range($ref).Value = $newResult
}

# --- Swap the "confirm" / "do not confirm" deletion steps between TC2 and TC3 ---
$ws.Range("B26").Value = "Lider de Pessoas nao confirma a exclusao do Periodo Avaliativo"
$ws.Range("D26").Value = "SYSTEM exibe a listagem dos Periodos Avaliativos com o Periodo Avaliativo nao excluido"

$ws.Range("B36").Value = "Lider de Pessoas confirma a exclusao do Periodo Avaliativo"
$ws.Range("D36").Value = "SYSTEM exibe a listagem dos Periodos Avaliativos sem o Periodo Avaliativo excluido"
